$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header names for the new columns S15:BB15 (MIMARKS air v4.0 fields)
$headers = @{
    "S" = 'barometric_press'
    "T" = 'biotic_relationship'
    "U" = 'carb_dioxide'
    "V" = 'carb_monoxide'
    "W" = 'chem_administration'
    "X" = 'elev'
    "Y" = 'extrachrom_elements'
    "Z" = 'humidity'
    "AA" = 'isolation_source'
    "AB" = 'methane'
    "AC" = 'misc_param'
    "AD" = 'organism_count'
    "AE" = 'oxy_stat_samp'
    "AF" = 'oxygen'
    "AG" = 'perturbation'
    "AH" = 'pollutants'
    "AI" = 'rel_to_oxygen'
    "AJ" = 'resp_part_matter'
    "AK" = 'samp_collect_device'
    "AL" = 'samp_mat_process'
    "AM" = 'samp_salinity'
    "AN" = 'samp_size'
    "AO" = 'samp_store_dur'
    "AP" = 'samp_store_loc'
    "AQ" = 'samp_store_temp'
    "AR" = 'samp_vol_we_dna_ext'
    "AS" = 'solar_irradiance'
    "AT" = 'source_material_id'
    "AU" = 'subspecf_gen_lin'
    "AV" = 'temperature'
    "AW" = 'trophic_level'
    "AX" = 'ventilation_rate'
    "AY" = 'ventilation_type'
    "AZ" = 'volatile_org_comp'
    "BA" = 'wind_direction'
    "BB" = 'wind_speed'
}

# Comment text (field definitions) for the same cells
$comments = @{
    "S" = 'force per unit area exerted against a surface by the weight of air above that surface'
    "T" = 'Free-living or from host (define relationship)'
    "U" = 'carbon dioxide (gas) amount or concentration at the time of sampling'
    "V" = 'carbon monoxide (gas) amount or concentration at the time of sampling'
    "W" = 'list of chemical compounds administered to the host or site where sampling occurred, and when (e.g. antibiotics, N fertilizer, air filter); can include multiple compounds. For Chemical Entities of Biological Interest ontology (CHEBI) (v1.72), please see http://bioportal.bioontology.org/visualize/44603'
    "X" = 'The elevation of the sampling site as measured by the vertical distance from mean sea level.'
    "Y" = 'Plasmids that have significance phenotypic consequence'
    "Z" = 'amount of water vapour in the air, at the time of sampling'
    "AA" = 'Describes the physical, environmental and/or local geographical source of the biological sample from which the sample was derived.'
    "AB" = 'methane (gas) amount or concentration at the time of sampling'
    "AC" = 'any other measurement performed or parameter collected, that is not listed here'
    "AD" = 'total count of any organism per gram or volume of sample, should include name of organism followed by count; can include multiple organism counts'
    "AE" = 'oxygenation status of sample'
    "AF" = 'oxygen (gas) amount or concentration at the time of sampling'
    "AG" = 'type of perturbation, e.g. chemical administration, physical disturbance, etc., coupled with time that perturbation occurred; can include multiple perturbation types'
    "AH" = 'pollutant types and, amount or concentrations measured at the time of sampling; can report multiple pollutants by entering numeric values preceded by name of pollutant'
    "AI" = 'Aerobic or anaerobic'
    "AJ" = 'concentration of substances that remain suspended in the air, and comprise mixtures of organic and inorganic substances (PM10 and PM2.5); can report multiple PM''s by entering numeric values preceded by name of PM'
    "AK" = 'Method or device employed for collecting sample'
    "AL" = 'Processing applied to the sample during or after isolation'
    "AM" = 'salinity of sample, i.e. measure of total salt concentration'
    "AN" = 'Amount or size of sample (volume, mass or area) that was collected'
    "AO" = 'duration for which sample was stored'
    "AP" = 'location at which sample was stored, usually name of a specific freezer/room'
    "AQ" = 'temperature at which sample was stored, e.g. -80'
    "AR" = 'volume (mL) or weight (g) of sample processed for DNA extraction'
    "AS" = 'the amount of solar energy that arrives at a specific area of a surface during a specific time interval'
    "AT" = 'unique identifier assigned to a material sample used for extracting nucleic acids, and subsequent sequencing. The identifier can refer either to the original material collected or to any derived sub-samples.'
    "AU" = 'Information about the genetic distinctness of the lineage (eg., biovar, serovar)'
    "AV" = 'temperature of the sample at time of sampling'
    "AW" = 'Feeding position in food chain (eg., chemolithotroph)'
    "AX" = 'ventilation rate of the system in the sampled premises'
    "AY" = 'ventilation system used in the sampled premises'
    "AZ" = 'concentration of carbon-based chemicals that easily evaporate at room temperature; can report multiple volatile organic compounds by entering numeric values preceded by name of compound'
    "BA" = 'wind direction is the direction from which a wind originates'
    "BB" = 'speed of wind measured at the time of sampling'
}

$cols = @("S", "T", "U", "V", "W", "X", "Y", "Z", "AA", "AB", "AC", "AD", "AE", "AF", "AG", "AH", "AI", "AJ", "AK", "AL", "AM", "AN", "AO", "AP", "AQ", "AR", "AS", "AT", "AU", "AV", "AW", "AX", "AY", "AZ", "BA", "BB")

foreach ($col in $cols) {
    $cell = $ws.Range($col + "15")
    $cell.Value = $headers[$col]
    # Match the yellow "optional field" header style already used on C15
    $ws.Range("C15").Copy()
    $cell.PasteSpecial(-4122)
    $cell.AddComment($comments[$col])
}

